# Update cryptos list prices and 1h volume percentages (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.550.68"
$ws.Range("E2").Value = "  -0.21%  "

$ws.Range("D3").Value = "1.728.83"
$ws.Range("E3").Value = "  -0.88%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'246.10"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("E7").Value = "  +0.54%  "

$ws.Range("D8").Value = "'0.2664"
$ws.Range("E8").Value = "  -1.05%  "

$ws.Range("D9").Value = "'0.06190"
$ws.Range("E9").Value = "  -1.00%  "

$ws.Range("D10").Value = "1.732.81"
$ws.Range("E10").Value = "  -0.46%  "

$ws.Range("D11").Value = "'0.07103"
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").Value = "'15.55"
$ws.Range("E12").Value = "  -1.53%  "

$ws.Range("D13").Value = "'0.6076"
$ws.Range("E13").Value = "  -1.76%  "

$ws.Range("D14").Value = "'4.544"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("D15").Value = "'77.16"
$ws.Range("E15").Value = "  -0.54%  "

$ws.Range("E16").Value = "  +0.09%  "

$ws.Range("D17").Value = "26.542.55"
$ws.Range("E17").Value = "  -0.22%  "

$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("D19").Value = "'0.000007220"
$ws.Range("E19").Value = "  +4.74%  "

$ws.Range("D20").Value = "'11.50"
$ws.Range("E20").Value = "  -1.80%  "

$ws.Range("D21").Value = "1.955.40"
$ws.Range("E21").Value = "  -0.60%  "

$ws.Range("D22").Value = "'4.504"
$ws.Range("E22").Value = "  -2.81%  "

$ws.Range("D23").Value = "'8.766"
$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("D24").Value = "'5.232"
$ws.Range("E24").Value = "  -2.19%  "

$ws.Range("D25").Value = "'137.41"
$ws.Range("E25").Value = "  +1.12%  "

$ws.Range("D26").Value = "'15.48"
$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("D27").Value = "'1.413"
$ws.Range("E27").Value = "  -1.66%  "

$ws.Range("D28").Value = "'1.773"
$ws.Range("E28").Value = "  -2.51%  "

$ws.Range("D29").Value = "'108.03"
$ws.Range("E29").Value = "  +0.65%  "

$ws.Range("D30").Value = "'3.957"
$ws.Range("E30").Value = "  -1.38%  "

$ws.Range("D31").Value = "'0.07999"
$ws.Range("E31").Value = "  +1.44%  "

$ws.Range("D32").Value = "'3.678"
$ws.Range("E32").Value = "  -2.01%  "

$ws.Range("D33").Value = "'0.04575"
$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").Value = "'1.000"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("E35").Value = "  +0.21%  "

$ws.Range("D36").Value = "'0.9985"
$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("D37").Value = "'0.6292"
$ws.Range("E37").Value = "  -1.84%  "

$ws.Range("D38").Value = "'0.8924"
$ws.Range("E38").Value = "  -5.47%  "

$ws.Range("D39").Value = "'2.005"
$ws.Range("E39").Value = "  +0.55%  "

$ws.Range("D40").Value = "'2.396"
$ws.Range("E40").Value = "  -1.78%  "

$ws.Range("E41").Value = "  +0.16%  "

$ws.Range("E42").Value = "  -0.49%  "

$ws.Range("D43").Value = "'101.65"
$ws.Range("E43").Value = "  -10.01%  "

$ws.Range("E44").Value = "  -6.37%  "

$ws.Range("D45").Value = "'0.3885"
$ws.Range("E45").Value = "  -0.76%  "

$ws.Range("D46").Value = "'7.042"
$ws.Range("E46").Value = "  +5.36%  "

$ws.Range("D47").Value = "'0.1183"
$ws.Range("E47").Value = "  -1.89%  "

$ws.Range("E48").Value = "  +1.21%  "

$ws.Range("D49").Value = "'7.906"
$ws.Range("E49").Value = "  -0.81%  "

$ws.Range("D50").Value = "'30.55"
$ws.Range("E50").Value = "  -0.96%  "

$ws.Range("D51").Value = "'1.258"
$ws.Range("E51").Value = "  -0.97%  "
